{"js": "// Replace the date line and each \"dividend\u00f7divisor=quotient, remainder\" cell\n// in the practice table with this commit's updated values.\nconst replacements = [\n  [\"2026-02-11 Wednesday\", \"2026-02-12 Thursday\"],\n  [\"614\u00f78=76, 6\", \"808\u00f77=115, 3\"],\n  [\"990\u00f72=495, 0\", \"515\u00f72=257, 1\"],\n  [\"404\u00f72=202, 0\", \"173\u00f77=24, 5\"],\n  [\"201\u00f76=33, 3\", \"418\u00f73=139, 1\"],\n  [\"709\u00f76=118, 1\", \"798\u00f79=88, 6\"],\n  [\"453\u00f72=226, 1\", \"547\u00f77=78, 1\"],\n  [\"431\u00f74=107, 3\", \"187\u00f77=26, 5\"],\n  [\"438\u00f78=54, 6\", \"351\u00f73=117, 0\"],\n  [\"443\u00f73=147, 2\", \"542\u00f75=108, 2\"],\n  [\"124\u00f79=13, 7\", \"221\u00f76=36, 5\"],\n  [\"811\u00f73=270, 1\", \"575\u00f74=143, 3\"],\n  [\"968\u00f72=484, 0\", \"356\u00f75=71, 1\"],\n  [\"894\u00f73=298, 0\", \"807\u00f72=403, 1\"],\n  [\"282\u00f78=35, 2\", \"752\u00f73=250, 2\"],\n  [\"584\u00f72=292, 0\", \"673\u00f72=336, 1\"],\n  [\"383\u00f75=76, 3\", \"761\u00f75=152, 1\"],\n  [\"742\u00f72=371, 0\", \"239\u00f76=39, 5\"],\n  [\"411\u00f77=58, 5\", \"790\u00f75=158, 0\"],\n  [\"534\u00f79=59, 3\", \"409\u00f78=51, 1\"],\n  [\"807\u00f76=134, 3\", \"988\u00f79=109, 7\"],\n  [\"444\u00f74=111, 0\", \"454\u00f73=151, 1\"],\n  [\"899\u00f78=112, 3\", \"277\u00f74=69, 1\"],\n  [\"922\u00f74=230, 2\", \"167\u00f79=18, 5\"],\n  [\"776\u00f72=388, 0\", \"946\u00f78=118, 2\"],\n  [\"629\u00f79=69, 8\", \"321\u00f73=107, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2026-02-11 Wednesday\", \"2026-02-12 Thursday\"),\n  @(\"614\u00f78=76, 6\", \"808\u00f77=115, 3\"),\n  @(\"990\u00f72=495, 0\", \"515\u00f72=257, 1\"),\n  @(\"404\u00f72=202, 0\", \"173\u00f77=24, 5\"),\n  @(\"201\u00f76=33, 3\", \"418\u00f73=139, 1\"),\n  @(\"709\u00f76=118, 1\", \"798\u00f79=88, 6\"),\n  @(\"453\u00f72=226, 1\", \"547\u00f77=78, 1\"),\n  @(\"431\u00f74=107, 3\", \"187\u00f77=26, 5\"),\n  @(\"438\u00f78=54, 6\", \"351\u00f73=117, 0\"),\n  @(\"443\u00f73=147, 2\", \"542\u00f75=108, 2\"),\n  @(\"124\u00f79=13, 7\", \"221\u00f76=36, 5\"),\n  @(\"811\u00f73=270, 1\", \"575\u00f74=143, 3\"),\n  @(\"968\u00f72=484, 0\", \"356\u00f75=71, 1\"),\n  @(\"894\u00f73=298, 0\", \"807\u00f72=403, 1\"),\n  @(\"282\u00f78=35, 2\", \"752\u00f73=250, 2\"),\n  @(\"584\u00f72=292, 0\", \"673\u00f72=336, 1\"),\n  @(\"383\u00f75=76, 3\", \"761\u00f75=152, 1\"),\n  @(\"742\u00f72=371, 0\", \"239\u00f76=39, 5\"),\n  @(\"411\u00f77=58, 5\", \"790\u00f75=158, 0\"),\n  @(\"534\u00f79=59, 3\", \"409\u00f78=51, 1\"),\n  @(\"807\u00f76=134, 3\", \"988\u00f79=109, 7\"),\n  @(\"444\u00f74=111, 0\", \"454\u00f73=151, 1\"),\n  @(\"899\u00f78=112, 3\", \"277\u00f74=69, 1\"),\n  @(\"922\u00f74=230, 2\", \"167\u00f79=18, 5\"),\n  @(\"776\u00f72=388, 0\", \"946\u00f78=118, 2\"),\n  @(\"629\u00f79=69, 8\", \"321\u00f73=107, 0\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    Write-Output \"NOT FOUND: $oldText\"\n  }\n}\n\n"}
